$wb = $excel.ActiveWorkbook

# --- Add the new "API_Controller" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "API_Controller"

$ws.Cells.Item(1,1).Value = "Trending  new panel"
$ws.Cells.Item(1,2).Value = "List<Article>"
$ws.Cells.Item(4,2).Value = "Article"
$ws.Cells.Item(1,3).Value = "Id"
$ws.Cells.Item(2,3).Value = "Channel"
$ws.Cells.Item(3,3).Value = "IconChannel"
$ws.Cells.Item(4,3).Value = "PubDate"
$ws.Cells.Item(5,3).Value = "Title"
$ws.Cells.Item(6,3).Value = "LikeNumber"
$ws.Cells.Item(7,3).Value = "DisLikeNumber"

$ws.Range("A1:A7").Merge()
$ws.Range("A1:A7").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A1:A7").VerticalAlignment = -4108    # xlCenter

$ws.Cells.Item(1,2).Font.Italic = $true

$ws.Cells.Item(4,2).HorizontalAlignment = -4131  # xlLeft

$ws.Rows.Item(1).RowHeight = 15.6

Write-Host "done"
